$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column D values can look like numbers (e.g. "1.008"); force them to stay
# text so they round-trip exactly like the source data, then restore the
# default "Normal" style so no stray number-format style is left behind.
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '26.805.68'
$ws.Range('E2').Value = '  +1.15%  '
Set-TextValue 'D3' '1.837.78'
$ws.Range('E3').Value = '  +1.56%  '
Set-TextValue 'D4' '1.008'
$ws.Range('E4').Value = '  +0.41%  '
Set-TextValue 'D5' '309.06'
$ws.Range('E5').Value = '  +1.34%  '
Set-TextValue 'D6' '1.007'
$ws.Range('E6').Value = '  +0.34%  '
Set-TextValue 'D7' '0.4701'
$ws.Range('E7').Value = '  +2.74%  '
$ws.Range('E8').Value = '  +1.81%  '
Set-TextValue 'D9' '0.07161'
$ws.Range('E9').Value = '  +0.91%  '
Set-TextValue 'D10' '0.9160'
$ws.Range('E10').Value = '  +2.05%  '
$ws.Range('B11').Value = 'WrappedEther'
$ws.Range('C11').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D11' '1.968.72'
$ws.Range('E11').Value = '  +8.95%  '
$ws.Range('B12').Value = 'Solana'
$ws.Range('C12').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue 'D12' '19.49'
$ws.Range('E12').Value = '  +1.07%  '
Set-TextValue 'D13' '0.07591'
$ws.Range('E13').Value = '  -1.78%  '
$ws.Range('E14').Value = '  +0.50%  '
Set-TextValue 'D15' '6.403'
$ws.Range('E15').Value = '  +1.86%  '
Set-TextValue 'D16' '87.58'
$ws.Range('E16').Value = '  +1.12%  '
Set-TextValue 'D17' '1.009'
$ws.Range('E17').Value = '  +0.37%  '
Set-TextValue 'D18' '0.000008610'
$ws.Range('E18').Value = '  +1.07%  '
Set-TextValue 'D19' '1.006'
$ws.Range('E19').Value = '  +0.36%  '
Set-TextValue 'D20' '26.848.98'
$ws.Range('E20').Value = '  +1.13%  '
Set-TextValue 'D21' '14.46'
$ws.Range('E21').Value = '  +2.57%  '
Set-TextValue 'D22' '5.000'
$ws.Range('E22').Value = '  +0.73%  '
$ws.Range('E23').Value = '  +0.56%  '
Set-TextValue 'D24' '1.924'
$ws.Range('E24').Value = '  +0.68%  '
Set-TextValue 'D25' '151.44'
$ws.Range('E25').Value = '  +0.08%  '
Set-TextValue 'D26' '18.10'
$ws.Range('E26').Value = '  +1.60%  '
Set-TextValue 'D27' '1.990'
$ws.Range('E27').Value = '  -0.87%  '
Set-TextValue 'D28' '113.76'
$ws.Range('E28').Value = '  +1.38%  '
Set-TextValue 'D29' '4.836'
$ws.Range('E29').Value = '  +0.75%  '
Set-TextValue 'D30' '0.08812'
$ws.Range('E30').Value = '  +0.87%  '
Set-TextValue 'D31' '3.304'
$ws.Range('E31').Value = '  +6.18%  '
Set-TextValue 'D32' '1.167'
$ws.Range('E32').Value = '  +5.04%  '
Set-TextValue 'D33' '0.7366'
$ws.Range('E33').Value = '  -0.41%  '
Set-TextValue 'D34' '4.472'
$ws.Range('E34').Value = '  +1.32%  '
Set-TextValue 'D35' '2.744'
$ws.Range('E35').Value = '  +1.07%  '
Set-TextValue 'D36' '1.084'
$ws.Range('E36').Value = '  +1.59%  '
Set-TextValue 'D37' '0.05249'
$ws.Range('E37').Value = '  +3.66%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D38' '2.978'
$ws.Range('E38').Value = '  +2.24%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D39' '0.01940'
$ws.Range('E39').Value = '  +0.83%  '
Set-TextValue 'D40' '0.5167'
$ws.Range('E40').Value = '  +2.47%  '
Set-TextValue 'D41' '6.876'
$ws.Range('E41').Value = '  +0.91%  '
Set-TextValue 'D42' '0.1506'
$ws.Range('E42').Value = '  +0.52%  '
Set-TextValue 'D43' '8.115'
$ws.Range('E43').Value = '  +2.06%  '
Set-TextValue 'D44' '10.43'
$ws.Range('E44').Value = '  +5.19%  '
Set-TextValue 'D45' '0.4667'
$ws.Range('E45').Value = '  +0.04%  '
Set-TextValue 'D46' '1.007'
$ws.Range('E46').Value = '  +0.41%  '
Set-TextValue 'D47' '101.37'
$ws.Range('E47').Value = '  +2.15%  '
Set-TextValue 'D48' '1.594'
$ws.Range('E48').Value = '  +2.34%  '
Set-TextValue 'D49' '65.49'
$ws.Range('E49').Value = '  +2.82%  '
Set-TextValue 'D50' '0.06025'
$ws.Range('E50').Value = '  +0.32%  '
Set-TextValue 'D51' '0.8851'
$ws.Range('E51').Value = '  +4.84%  '
